# "fixed a bug of Bard": remove the obsolete 3rd trigger group
# (comment "#第二轮先切歌,再延后团辅" + its Cond:AfterBattleStart / Action:CastAbility
# rows) and shift the following groups up. Also correct the burst-alignment
# trigger's timing value now that it moved into the vacated slot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "第二轮先切歌,再延后团辅" group: its header row (13) and the
# two data rows underneath it (14:15) that referenced Action:CastAbility.
$ws.Rows("13:15").Delete()

# The group that used to sit right after it (comment row 16 -> "#回场中开爆发,对齐150秒的团辅")
# now starts the block at row 13; its Cond:AfterBattleStart value changes from 145 to 140.
$ws.Range("F14").Value = 140

# Leave the selection where the author left off editing.
$ws.Range("F15").Select()
